# Apply cryptos list update (price/volume refresh + two name/row swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '37.294.68'
Set-TextValue 'E2' '  +1.72%  '
Set-TextValue 'D3' '2.061.41'
Set-TextValue 'E3' '  +1.24%  '
Set-TextValue 'E4' '  -0.01%  '
Set-TextValue 'D5' '232.44'
Set-TextValue 'E5' '  -0.13%  '
Set-TextValue 'E6' '  +2.81%  '
Set-TextValue 'D8' '57.14'
Set-TextValue 'E8' '  +3.28%  '
Set-TextValue 'E9' '  +3.07%  '
Set-TextValue 'D10' '58.06'
Set-TextValue 'E10' '  +0.84%  '
Set-TextValue 'D11' '0.0761'
Set-TextValue 'E11' '  +0.79%  '
Set-TextValue 'E12' '  +1.28%  '
Set-TextValue 'D13' '2.364.53'
Set-TextValue 'E13' '  +1.29%  '
Set-TextValue 'D14' '14.54'
Set-TextValue 'E14' '  +1.96%  '
Set-TextValue 'D15' '20.86'
Set-TextValue 'E15' '  +3.93%  '
Set-TextValue 'D16' '0.779'
Set-TextValue 'E16' '  +2.06%  '
Set-TextValue 'E17' '  -0.10%  '
Set-TextValue 'D18' '2.059.93'
Set-TextValue 'E18' '  +1.53%  '
Set-TextValue 'D19' '37.209.61'
Set-TextValue 'E19' '  +1.19%  '
Set-TextValue 'D20' '6.31'
Set-TextValue 'E20' '  +7.15%  '
Set-TextValue 'D21' '69.23'
Set-TextValue 'E21' '  +2.08%  '
Set-TextValue 'D22' '0.0₃0811'
Set-TextValue 'E22' '  +1.26%  '
Set-TextValue 'D23' '225.96'
Set-TextValue 'E23' '  +2.28%  '
Set-TextValue 'E24' '  +0.03%  '
Set-TextValue 'E25' '  +0.64%  '
Set-TextValue 'D26' '2.39'
Set-TextValue 'E26' '  +0.50%  '
Set-TextValue 'D27' '166.08'
Set-TextValue 'E27' '  +1.95%  '
Set-TextValue 'E28' '  +8.40%  '
Set-TextValue 'E29' '  +0.23%  '
Set-TextValue 'D30' '0.128'
Set-TextValue 'E30' '  +0.94%  '
Set-TextValue 'E31' '  +0.51%  '
Set-TextValue 'D32' '0.118'
Set-TextValue 'E32' '  +0.17%  '
Set-TextValue 'D33' '4.47'
Set-TextValue 'E33' '  +2.26%  '
Set-TextValue 'B34' 'Hedera'
Set-TextValue 'C34' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D34' '0.0619'
Set-TextValue 'E34' '  +1.59%  '
Set-TextValue 'B35' 'InternetComputer(DFINITY)'
Set-TextValue 'C35' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D35' '4.62'
Set-TextValue 'E35' '  +7.74%  '
Set-TextValue 'E36' '  +1.66%  '
Set-TextValue 'E37' '  +0.06%  '
Set-TextValue 'D38' '1.75'
Set-TextValue 'E38' '  -0.25%  '
Set-TextValue 'D39' '3.27'
Set-TextValue 'E39' '  +1.76%  '
Set-TextValue 'D40' '5.69'
Set-TextValue 'E40' '  -1.59%  '
Set-TextValue 'E41' '  +0.30%  '
Set-TextValue 'E42' '  +2.87%  '
Set-TextValue 'D43' '1.466.61'
Set-TextValue 'E43' '  -1.15%  '
Set-TextValue 'D44' '4.33'
Set-TextValue 'E44' '  -1.42%  '
Set-TextValue 'B45' 'VeChain'
Set-TextValue 'C45' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D45' '0.0213'
Set-TextValue 'E45' '  +4.13%  '
Set-TextValue 'D46' '0.0932'
Set-TextValue 'E46' '  -1.59%  '
Set-TextValue 'B47' 'TrustWalletToken'
Set-TextValue 'C47' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D47' '1.17'
Set-TextValue 'E47' '  +5.40%  '
Set-TextValue 'E48' '  +1.06%  '
Set-TextValue 'D49' '15.12'
Set-TextValue 'E49' '  -3.23%  '
Set-TextValue 'D50' '7.17'
Set-TextValue 'E50' '  +2.75%  '
Set-TextValue 'E51' '  +2.02%  '
